# Update "想去人数" (interested-people count) figures in the 展览 and 全部类型
# sheets to reflect newly generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1840
$ws1.Range("F13").Value = 535
$ws1.Range("F14").Value = 526
$ws1.Range("F20").Value = 1601
$ws1.Range("F32").Value = 3812
$ws1.Range("F33").Value = 760
$ws1.Range("F35").Value = 775
$ws1.Range("F37").Value = 1812

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1840
$ws4.Range("F13").Value = 535
$ws4.Range("F14").Value = 526
$ws4.Range("F21").Value = 1601
$ws4.Range("F33").Value = 3812
$ws4.Range("F35").Value = 760
$ws4.Range("F37").Value = 776
$ws4.Range("F39").Value = 1812
